$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.244.03'
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.841.69'
$ws.Range("E3").Value = '  +1.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '361.08'
$ws.Range("E5").Value = '  +5.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.30'
$ws.Range("E6").Value = '  -2.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.577'
$ws.Range("E7").Value = '  +4.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  +4.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.61'
$ws.Range("E10").Value = '  -1.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.132'
$ws.Range("E12").Value = '  +1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.05'
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.82'
$ws.Range("E14").Value = '  +2.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.287.39'
$ws.Range("E15").Value = '  +1.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.835.51'
$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.910'
$ws.Range("E17").Value = '  +2.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.125.10'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.62'
$ws.Range("E19").Value = '  +9.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.15'
$ws.Range("E20").Value = '  -2.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.55'
$ws.Range("E21").Value = '  +1.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0997'
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.47'
$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.51'
$ws.Range("E24").Value = '  -3.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.84'
$ws.Range("E25").Value = '  +2.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.18'
$ws.Range("E26").Value = '  +1.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.43'
$ws.Range("E28").Value = '  +1.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +1.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '54.50'
$ws.Range("E30").Value = '  +7.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0480'
$ws.Range("E31").Value = '  +27.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.141'
$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.75'
$ws.Range("E33").Value = '  -0.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("E34").Value = '  +2.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.49'
$ws.Range("E35").Value = '  +10.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0846'
$ws.Range("E36").Value = '  +2.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("E38").Value = '  +0.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.08'
$ws.Range("E39").Value = '  -2.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.41'
$ws.Range("E40").Value = '  -3.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.97'
$ws.Range("E41").Value = '  +2.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").Value = '  +1.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '128.00'
$ws.Range("E43").Value = '  +2.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("E44").Value = '  -7.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.28'
$ws.Range("E45").Value = '  -1.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.43'
$ws.Range("E46").Value = '  +2.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.115.05'
$ws.Range("E47").Value = '  +0.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.26'
$ws.Range("E48").Value = '  +1.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("E49").Value = '  +11.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.88'
$ws.Range("E50").Value = '  +5.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.02'
$ws.Range("E51").Value = '  +2.91%  '
